$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.063063333333333
$ws.Range("H2").Value = 6.189190000000001
$ws.Range("I2").Value = 0.1875199417503197
$ws.Range("J2").Value = 0.1875199417503197
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 104.9303403108956
$ws.Range("R2").Value = 944.3730627980601
$ws.Range("S2").Value = 0.06440024678587346
$ws.Range("T2").Value = 0.06440024678587347
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.063063333333333
$ws.Range("H3").Value = 6.189190000000001
$ws.Range("I3").Value = 0.1875199417503197
$ws.Range("J3").Value = 0.1875199417503197
$ws.Range("M3").Value = 43.683024
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("Q3").Value = 90.12084510352
$ws.Range("R3").Value = 811.0876059316801
$ws.Range("S3").Value = 0.05531102489539452
$ws.Range("T3").Value = 0.05531102489539454
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.063063333333333
$ws.Range("H4").Value = 6.189190000000001
$ws.Range("I4").Value = 0.1875199417503197
$ws.Range("J4").Value = 0.1875199417503197
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 75.59807925202556
$ws.Range("R4").Value = 680.3827132682301
$ws.Range("S4").Value = 0.046397781098809
$ws.Range("T4").Value = 0.04639778109880902
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.063063333333333
$ws.Range("H5").Value = 6.189190000000001
$ws.Range("I5").Value = 0.1875199417503197
$ws.Range("J5").Value = 0.1875199417503197
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 34.88576485547222
$ws.Range("R5").Value = 313.97188369925
$ws.Range("S5").Value = 0.02141088897024266
$ws.Range("T5").Value = 0.02141088897024267
$ws.Range("I6").Value = 0.5238509469163369
$ws.Range("J6").Value = 0.5238509469163369
$ws.Range("M6").Value = 50.86142466666666
$ws.Range("N6").Value = 152.584274
$ws.Range("O6").Value = 0.3434314568613803
$ws.Range("P6").Value = 0.3434314568613804
$ws.Range("Q6").Value = 293.1307338251262
$ws.Range("R6").Value = 2638.176604426135
$ws.Range("S6").Value = 0.1799068938776912
$ws.Range("T6").Value = 0.1799068938776912
$ws.Range("I7").Value = 0.5238509469163369
$ws.Range("J7").Value = 0.5238509469163369
$ws.Range("M7").Value = 43.683024
$ws.Range("O7").Value = 0.294960761928139
$ws.Range("P7").Value = 0.294960761928139
$ws.Range("S7").Value = 0.1545154744392198
$ws.Range("T7").Value = 0.1545154744392198
$ws.Range("I8").Value = 0.5238509469163369
$ws.Range("J8").Value = 0.5238509469163369
$ws.Range("M8").Value = 36.64360566666667
$ws.Range("N8").Value = 109.930817
$ws.Range("O8").Value = 0.2474285170192034
$ws.Range("P8").Value = 0.2474285170192035
$ws.Range("Q8").Value = 211.1888742689542
$ws.Range("R8").Value = 1900.699868420588
$ws.Range("S8").Value = 0.1296156629346147
$ws.Range("T8").Value = 0.1296156629346147
$ws.Range("I9").Value = 0.5238509469163369
$ws.Range("J9").Value = 0.5238509469163369
$ws.Range("M9").Value = 16.90969166666667
$ws.Range("N9").Value = 50.729075
$ws.Range("O9").Value = 0.1141792641912772
$ws.Range("P9").Value = 0.1141792641912772
$ws.Range("Q9").Value = 97.45598672258889
$ws.Range("R9").Value = 877.1038805032999
$ws.Range("S9").Value = 0.05981291566481117
$ws.Range("T9").Value = 0.05981291566481117
$ws.Range("G10").Value = 2.101774
$ws.Range("H10").Value = 6.305322
$ws.Range("I10").Value = 0.1910385065181404
$ws.Range("J10").Value = 0.1910385065181404
$ws.Range("M10").Value = 50.86142466666666
$ws.Range("N10").Value = 152.584274
$ws.Range("O10").Value = 0.3434314568613803
$ws.Range("P10").Value = 0.3434314568613804
$ws.Range("Q10").Value = 106.8992199673587
$ws.Range("R10").Value = 962.092979706228
$ws.Range("S10").Value = 0.06560863261014724
$ws.Range("T10").Value = 0.06560863261014725
$ws.Range("G11").Value = 2.101774
$ws.Range("H11").Value = 6.305322
$ws.Range("I11").Value = 0.1910385065181404
$ws.Range("J11").Value = 0.1910385065181404
$ws.Range("M11").Value = 43.683024
$ws.Range("O11").Value = 0.294960761928139
$ws.Range("P11").Value = 0.294960761928139
$ws.Range("Q11").Value = 91.81184408457601
$ws.Range("R11").Value = 826.306596761184
$ws.Range("S11").Value = 0.05634886344020442
$ws.Range("T11").Value = 0.05634886344020443
$ws.Range("G12").Value = 2.101774
$ws.Range("H12").Value = 6.305322
$ws.Range("I12").Value = 0.1910385065181404
$ws.Range("J12").Value = 0.1910385065181404
$ws.Range("M12").Value = 36.64360566666667
$ws.Range("N12").Value = 109.930817
$ws.Range("O12").Value = 0.2474285170192034
$ws.Range("P12").Value = 0.2474285170192035
$ws.Range("Q12").Value = 77.01657765645267
$ws.Range("R12").Value = 693.149198908074
$ws.Range("S12").Value = 0.0472683743613469
$ws.Range("T12").Value = 0.04726837436134691
$ws.Range("G13").Value = 2.101774
$ws.Range("H13").Value = 6.305322
$ws.Range("I13").Value = 0.1910385065181404
$ws.Range("J13").Value = 0.1910385065181404
$ws.Range("M13").Value = 16.90969166666667
$ws.Range("N13").Value = 50.729075
$ws.Range("O13").Value = 0.1141792641912772
$ws.Range("P13").Value = 0.1141792641912772
$ws.Range("Q13").Value = 35.54035029301667
$ws.Range("R13").Value = 319.8631526371501
$ws.Range("S13").Value = 0.02181263610644178
$ws.Range("T13").Value = 0.02181263610644178
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.073675666666667
$ws.Range("H14").Value = 3.221027
$ws.Range("I14").Value = 0.0975906048152031
$ws.Range("J14").Value = 0.09759060481520311
$ws.Range("M14").Value = 50.86142466666666
$ws.Range("N14").Value = 152.584274
$ws.Range("O14").Value = 0.3434314568613803
$ws.Range("P14").Value = 0.3434314568613804
$ws.Range("Q14").Value = 54.60867403659978
$ws.Range("R14").Value = 491.478066329398
$ws.Range("S14").Value = 0.03351568358766844
$ws.Range("T14").Value = 0.03351568358766845
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.073675666666667
$ws.Range("H15").Value = 3.221027
$ws.Range("I15").Value = 0.0975906048152031
$ws.Range("J15").Value = 0.09759060481520311
$ws.Range("M15").Value = 43.683024
$ws.Range("O15").Value = 0.294960761928139
$ws.Range("P15").Value = 0.294960761928139
$ws.Range("Q15").Value = 46.901399915216
$ws.Range("R15").Value = 422.112599236944
$ws.Range("S15").Value = 0.02878539915332021
$ws.Range("T15").Value = 0.02878539915332022
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.073675666666667
$ws.Range("H16").Value = 3.221027
$ws.Range("I16").Value = 0.0975906048152031
$ws.Range("J16").Value = 0.09759060481520311
$ws.Range("M16").Value = 36.64360566666667
$ws.Range("N16").Value = 109.930817
$ws.Range("O16").Value = 0.2474285170192034
$ws.Range("P16").Value = 0.2474285170192035
$ws.Range("Q16").Value = 39.34334774322878
$ws.Range("R16").Value = 354.090129689059
$ws.Range("S16").Value = 0.02414669862443284
$ws.Range("T16").Value = 0.02414669862443284
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.073675666666667
$ws.Range("H17").Value = 3.221027
$ws.Range("I17").Value = 0.0975906048152031
$ws.Range("J17").Value = 0.09759060481520311
$ws.Range("M17").Value = 16.90969166666667
$ws.Range("N17").Value = 50.729075
$ws.Range("O17").Value = 0.1141792641912772
$ws.Range("P17").Value = 0.1141792641912772
$ws.Range("Q17").Value = 18.15552447333611
$ws.Range("R17").Value = 163.399720260025
$ws.Range("S17").Value = 0.0111428234497816
$ws.Range("T17").Value = 0.01114282344978161